# Add HURM-009 "Actualización de Datos del Empleado" user-story block,
# right after the existing HURM-008 block (after the paragraph that ends
# with "Para: Verificar su información o actualizarla si es necesario").
#
# Final layout (matching the commit's diff) becomes:
#   ... Para: Verificar ... (existing, unchanged)
#   <empty paragraph>                                    (new)
#   HURM-009 Actualización de Datos del Empleado          (new)
#   Como: Empleado del área / Quiero: ... / Para: ...     (new)
#   <empty paragraph>                                     (new)
#   <the two pre-existing empty paragraphs>               (unchanged)
#   <Arial/sectPr paragraph>                              (unchanged)

$d = $word.ActiveDocument

# Locate the end of the HURM-008 block via its last line of text.
$r = $d.Content
$found = $r.Find.Execute(
    "Para: Verificar su información o actualizarla si es necesario",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph 'Para: Verificar su información o actualizarla si es necesario'"
}

# Collapse the found range to its end (right after the last character of
# that paragraph's text, i.e. right before its own paragraph mark) and
# rebuild a fresh Range from those coordinates -- reusing the Find-bound
# Range object directly for InsertXML can make the engine swallow the
# anchor paragraph, so we hand it a brand-new Range instead.
$r.Collapse(0)
$insertPos = $r.Start
$target = $d.Range($insertPos, $insertPos)

$huPara = '<w:p>' +
    '<w:r><w:t>HU</w:t></w:r>' +
    '<w:r><w:t>RM</w:t></w:r>' +
    '<w:r><w:t>-00</w:t></w:r>' +
    '<w:r><w:t>9</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Actualización de Datos del Empleado</w:t></w:r>' +
    '</w:p>'

$comoPara = '<w:p>' +
    '<w:r><w:t xml:space="preserve">Como: </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Empleado </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">del </w:t></w:r>' +
    '<w:r><w:t>área</w:t></w:r>' +
    '<w:r><w:br/><w:t>Quiero: Modificar o actualizar los datos de un empleado existente</w:t></w:r>' +
    '<w:r><w:br/><w:t>Para: Mantener la información del personal siempre actualizada</w:t></w:r>' +
    '</w:p>'

$xmlSnippet = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p/>' +
    $huPara +
    $comoPara +
    '<w:p/>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

[void]$target.InsertXML($xmlSnippet)
